$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/%) must stay text cells exactly as scraped,
# even though their contents look numeric/percentage. Force text format first
# so Excel does not silently convert them to numbers, then clear the transient
# number-format override so the cell keeps the workbook default style.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '313.23'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '2.06%'
$ws.Range('E2').ClearFormats()
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '39.94'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '-2.80%'
$ws.Range('E3').ClearFormats()
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.180'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '-0.99%'
$ws.Range('E4').ClearFormats()
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.07574'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '-1.11%'
$ws.Range('E5').ClearFormats()
$ws.Range('B6').Value = 'FTXToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.681'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '3.44%'
$ws.Range('E6').ClearFormats()
$ws.Range('B7').Value = 'MXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9255'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '0.96%'
$ws.Range('E7').ClearFormats()
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '2.424'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '-0.82%'
$ws.Range('E8').ClearFormats()
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.1201'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '-4.43%'
$ws.Range('E9').ClearFormats()
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1816'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-0.53%'
$ws.Range('E10').ClearFormats()
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.09015'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '-2.04%'
$ws.Range('E11').ClearFormats()
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.04164'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-2.76%'
$ws.Range('E12').ClearFormats()
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.1052'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '0.12%'
$ws.Range('E13').ClearFormats()
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.001277'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '0.75%'
$ws.Range('E14').ClearFormats()
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.005835'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '-0.04%'
$ws.Range('E15').ClearFormats()
$ws.Range('B16').Value = 'UpBots'
$ws.Range('C16').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.007522'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '0.18%'
$ws.Range('E16').ClearFormats()
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.350'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '-0.17%'
$ws.Range('E17').ClearFormats()
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '4.350'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '0.62%'
$ws.Range('E18').ClearFormats()
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3356'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '0.63%'
$ws.Range('E19').ClearFormats()
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.602'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '6.23%'
$ws.Range('E20').ClearFormats()
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '-2.39%'
$ws.Range('E21').ClearFormats()
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.2809'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-2.98%'
$ws.Range('E22').ClearFormats()
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04016'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '-1.34%'
$ws.Range('E23').ClearFormats()
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.001267'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '0.30%'
$ws.Range('E24').ClearFormats()
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.003978'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '-4.48%'
$ws.Range('E25').ClearFormats()
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0001270'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '-0.41%'
$ws.Range('E26').ClearFormats()
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '-1.46%'
$ws.Range('E38').ClearFormats()
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05151'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '-2.57%'
$ws.Range('E39').ClearFormats()
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.007701'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '-1.90%'
$ws.Range('E40').ClearFormats()
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.1300'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '-1.08%'
$ws.Range('E41').ClearFormats()
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.007618'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '11.49%'
$ws.Range('E42').ClearFormats()
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '51.39%'
$ws.Range('E43').ClearFormats()
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.008204'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '5.32%'
$ws.Range('E44').ClearFormats()
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.3109'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '1.84%'
$ws.Range('E45').ClearFormats()
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00006585'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '-1.85%'
$ws.Range('E46').ClearFormats()
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-0.33%'
$ws.Range('E47').ClearFormats()
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.2697'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '58.75%'
$ws.Range('E48').ClearFormats()
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.004199'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '2.41%'
$ws.Range('E49').ClearFormats()
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.00002099'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '-0.33%'
$ws.Range('E50').ClearFormats()
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0001999'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '-0.33%'
$ws.Range('E51').ClearFormats()
